$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 41

# Columns A (date) and D (week number) look numeric/date-like to Excel's
# automatic type detection, so force them to be stored as text (matching
# the rest of the column) using a leading apostrophe, just like the
# existing rows above which are plain text values.
$ws.Cells.Item($row, 1).Value = "'2023-06-12"
$ws.Cells.Item($row, 2).Value = "17:32:42"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "'24"
$ws.Cells.Item($row, 5).Value = 121403
$ws.Cells.Item($row, 6).Value = 135097
$ws.Cells.Item($row, 7).Value = 161161
$ws.Cells.Item($row, 8).Value = 132012
$ws.Cells.Item($row, 9).Value = 176428
$ws.Cells.Item($row, 10).Value = 113879
$ws.Cells.Item($row, 11).Value = 202052
$ws.Cells.Item($row, 12).Value = 222715
$ws.Cells.Item($row, 13).Value = 174058
$ws.Cells.Item($row, 14).Value = 100735
$ws.Cells.Item($row, 15).Value = 38807
$ws.Cells.Item($row, 16).Value = 34099
$ws.Cells.Item($row, 17).Value = 51309
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36828
$ws.Cells.Item($row, 20).Value = -1
